$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (U2021-0003 / ADAM HAIKAL) - entry no longer present after edit
$ws.Rows.Item(3).Delete()

# Full target dataset for rows 2-14 (user_id, name, reason)
$data = @(
    @("U2020-0058", "MUHAMMAD ZIDNI", "There are more than 1 data"),
    @("U2021-0004", "ADIB M P", "Data doesnt exist"),
    @("U2021-0007", "AGUS NURYADI", "There are more than 1 data"),
    @("U2021-0175", "NURUL FITRI", "There are more than 1 data"),
    @("U2021-0188", "RAIHAN ABDAN", "There are more than 1 data"),
    @("U2022-0021", "ELSYE", "There are more than 1 data"),
    @("U2022-0027", "HANI MAULIDIYAH", "Data doesnt exist"),
    @("U2022-0028", "INGE", "There are more than 1 data"),
    @("U2022-0032", "KENNETH LAURENTIUS", "There are more than 1 data"),
    @("U2022-0033", "KHARISMA YOGI PRATIWI", "Data doesnt exist"),
    @("U2022-0034", "KHUSNUL CHOTIMAH", "There are more than 1 data"),
    @("U2022-0036", "MEYSARA", "There are more than 1 data"),
    @("U2022-0038", "MISBAHUDDIN", "There are more than 1 data")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
